$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.919583201408386
$ws.Range("B1").Value = 2.753126382827759
$ws.Range("C1").Value = 3.290869474411011
$ws.Range("D1").Value = 1.080532431602478
$ws.Range("E1").Value = 0.6961435079574585
